$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, matching the bold/bordered header style
# used by the other header cells (B1:G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value2 = "Save"

# Find the last used data row from column A (dates run from row 2 to 50).
$lastRow = $ws.Cells.Item($ws.Rows.Count, "A").End(-4162).Row

# Populate H2:H<lastRow> with a save flag: 1 when the row's sum (column G)
# reaches a full game's worth (>= 9), otherwise 0.
for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ($g -ge 9) {
        $ws.Cells.Item($r, 8).Value2 = 1
    } else {
        $ws.Cells.Item($r, 8).Value2 = 0
    }
}
